$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.000.94'
$ws.Range("E2").Value = '  -3.80%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.172.91'
$ws.Range("E3").Value = '  -4.04%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.71'
$ws.Range("E5").Value = '  -5.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.43'
$ws.Range("E6").Value = '  -6.40%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.175.64'
$ws.Range("E8").Value = '  -4.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.453'
$ws.Range("E9").Value = '  -6.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.36'
$ws.Range("E10").Value = '  -5.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.112'
$ws.Range("E11").Value = '  -6.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.392'
$ws.Range("E12").Value = '  -3.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.713.83'
$ws.Range("E13").Value = '  -4.27%  '
$ws.Range("E14").Value = '  -1.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.94'
$ws.Range("E15").Value = '  -4.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.162.68'
$ws.Range("E16").Value = '  -4.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '57.943.65'
$ws.Range("E17").Value = '  -3.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000154'
$ws.Range("E18").Value = '  -7.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.84'
$ws.Range("E19").Value = '  -5.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.13'
$ws.Range("E20").Value = '  -8.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.12'
$ws.Range("E21").Value = '  -5.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '349.08'
$ws.Range("E22").Value = '  -6.58%  '
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.82'
$ws.Range("E24").Value = '  -5.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.514'
$ws.Range("E25").Value = '  -6.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.292.87'
$ws.Range("E26").Value = '  -5.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0968'
$ws.Range("E27").Value = '  -8.16%  '
$ws.Range("E28").Value = '  -3.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.92'
$ws.Range("E30").Value = '  -4.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("E32").Value = '  -7.84%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.97'
$ws.Range("E33").Value = '  -8.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.84'
$ws.Range("E34").Value = '  -3.29%  '
$ws.Range("E35").Value = '  -4.75%  '
$ws.Range("E36").Value = '  -4.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.61'
$ws.Range("E37").Value = '  -3.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.29'
$ws.Range("E38").Value = '  -7.05%  '
$ws.Range("E39").Value = '  -7.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.00'
$ws.Range("E40").Value = '  -6.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0700'
$ws.Range("E41").Value = '  -5.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.199.90'
$ws.Range("E42").Value = '  -4.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.74'
$ws.Range("E43").Value = '  -2.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.701'
$ws.Range("E44").Value = '  -6.79%  '
$ws.Range("E45").Value = '  -2.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.97'
$ws.Range("E46").Value = '  -5.71%  '
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("E48").Value = '  -7.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.275.25'
$ws.Range("E49").Value = '  -4.45%  '
$ws.Range("E50").Value = '  -5.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.69'
$ws.Range("E51").Value = '  -4.60%  '
